# Update column F ("dSF") values that were repulled/recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    6  = 0
    16 = -1
    17 = -5
    27 = 3
    28 = -9
    29 = 3
    30 = -3
    33 = -1
    34 = 5
    35 = 0
    36 = -1
    41 = -1
    44 = 6
    45 = -1
    46 = -1
    47 = 0
    49 = 3
    51 = 1
    53 = -2
    54 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
